# Apply "add pre trig window, update doc" changes to the firmware register map.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header: update date and clarify register-space note ---
$ws.Range("A4").Value = (Get-Date -Year 2017 -Month 7 -Day 1 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("C4").Value = "current 32 bit register space: 8 bit address + 24 bit data (i.e. 128 24-bit registers)"

# --- status_0 register: note that nothing is defined there yet ---
$ws.Range("E10").Value = "nothing defined here yet"

# --- chip_id registers: clarify read_only wording ---
$ws.Range("C11").Value = "read_only - board DNA (FPGA silicon-specific ID)"
$ws.Range("C12").Value = "read_only - board DNA (FPGA silicon-specific ID)"
$ws.Range("C13").Value = "read_only - board DNA (FPGA silicon-specific ID)"

# --- fill in functionality column for the unused/reserved register rows ---
$ws.Range("C14").Value = "read_only"
$ws.Range("C15").Value = "read_only"
$ws.Range("C16").Value = "read_only"
$ws.Range("C17").Value = "read_only"
$ws.Range("C18").Value = "read_only"
$ws.Range("C19").Value = "read_only"
$ws.Range("C20").Value = "read_only"
$ws.Range("C21").Value = "read_only"
$ws.Range("C22").Value = "read_only"
$ws.Range("C23").Value = "read_only"
$ws.Range("C24").Value = "read_only"
$ws.Range("C25").Value = "read_only"
$ws.Range("C26").Value = "read_only"
$ws.Range("C27").Value = "read_only"
$ws.Range("C28").Value = "read_only"
$ws.Range("C29").Value = "read_only"
$ws.Range("C30").Value = "read_only"
$ws.Range("C31").Value = "read_only"
$ws.Range("C32").Value = "read_only"
$ws.Range("C33").Value = "read_only"
$ws.Range("C34").Value = "read_only"
$ws.Range("C35").Value = "read_only"
$ws.Range("C36").Value = "read_only"
$ws.Range("C37").Value = "read_only"
$ws.Range("C38").Value = "read_only"

# --- clock select note: drop the "EVENTUALLY SWITCH..." reminder ---
$ws.Range("E39").Value = "lower bit only -- "

# --- readout register: BBB interface (was MCU) ---
$ws.Range("E78").Value = "for BBB interface, write this register to initiate readout"

# --- new register: rdout->pretrigger window ---
$ws.Range("B83").Value = "rdout->pretrigger window"
$ws.Range("C83").Value = "lower three bits used, set value to 0 thru 5. Pretrig window set to value*8*10.66ns"
$ws.Range("D83").Value = "x000000"

# --- thresholds: note default max value ---
$ws.Range("E93").Value = "start thresholds off at max"

# --- reset the view: select C6 and scroll back to top ---
$ws.Range("C6").Select()
